$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 2016 data set, ordered to match target row positions (rows 2-5)
$data = @(
    @(2016, 293105,  "University Hospital",                       "BEXAR",   "TX", 6940900,   1144160290, 1404986559, 1330472170, 2735458729, 136750936, 503605279,  0.4172113640553589, 503605279,  640356215,  267163889.9414766, 210109845.3970696, 210109845.3970696),
    @(2016, 1130950, "Parkland Memorial Hospital",                 "DALLAS",  "TX", 131483015, 1631498659, 2315400264, 3504218877, 5819619141, 607674875, 2110198804, 0.2741506726371847, 2110198804, $null,      $null,              578512421.5147827, 578512421.5147827),
    @(2016, 1411240, "University Medical Center of El Paso",       "EL PASO", "TX", 242070017, 608801157,  492781458,  514102315,  1006883773, 119187764, 225165406,  0.4874489047348982, 224765380,  344353170,  167854575.5584902, 109756630.5388887, 109561638.3033232),
    @(2016, 2015024, "Harris Health System Ben Taub Hospital",     "HARRIS",  "TX", 28065601,  1301032626, 1009831942, 2094350203, 3104182145, 115154658, 1547833906, 0.4153670882711843, 1547833905, 1662988564, 690750717.656958,  642919262.662634,  642919262.2472669)
)

$rowIndex = 2
foreach ($rec in $data) {
    for ($col = 1; $col -le 18; $col++) {
        $cell = $ws.Cells.Item($rowIndex, $col)
        $val = $rec[$col - 1]
        if ($null -eq $val) {
            $cell.Value = $null
        } else {
            $cell.Value = $val
        }
    }
    $rowIndex++
}
